$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 0.000625
$ws.Range("K2").Value = 4425
$ws.Range("L2").Value = 0.00885
